$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2900107.2
$ws.Range("I113").Value = 5129904.5
$ws.Range("J113").Value = 1370.6
$ws.Range("K113").Value = 5129904.5
$ws.Range("L113").Value = 1370.6
$ws.Range("M113").Value = -5126650.5
$ws.Range("N113").Value = -7878.6
$ws.Range("H137").Value = 898.21875
$ws.Range("I137").Value = 924.8077
$ws.Range("J137").Value = 783
$ws.Range("K137").Value = 2774.4231
$ws.Range("L137").Value = 2349
$ws.Range("M137").Value = -224.4231
$ws.Range("N137").Value = -7449

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 8731.691999999999
$ws.Range("I45").Value = 10799.8
$ws.Range("K45").Value = 10799.8
$ws.Range("M45").Value = -10422.8
$ws.Range("H97").Value = 654.96
$ws.Range("I97").Value = 407.1
$ws.Range("J97").Value = 1646.4
$ws.Range("K97").Value = 407.1
$ws.Range("L97").Value = 1646.4
$ws.Range("M97").Value = 88.89999999999998
$ws.Range("N97").Value = -2638.4
$ws.Range("H122").Value = 2071
$ws.Range("I122").Value = 2107.8948
$ws.Range("J122").Value = 1370
$ws.Range("K122").Value = 6323.6844
$ws.Range("L122").Value = 4110
$ws.Range("M122").Value = -3873.6844
$ws.Range("N122").Value = -9010

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 13268
$ws.Range("I20").Value = 1255.4615
$ws.Range("K20").Value = 1255.4615
$ws.Range("M20").Value = -1008.4615
$ws.Range("H99").Value = 1337.8422
$ws.Range("I99").Value = 913
$ws.Range("J99").Value = 1809.8889
$ws.Range("K99").Value = 913
$ws.Range("L99").Value = 1809.8889
$ws.Range("M99").Value = 585
$ws.Range("N99").Value = -4805.8889
$ws.Range("H105").Value = 12061.143
$ws.Range("I105").Value = 16659.715
$ws.Range("K105").Value = 16659.715
$ws.Range("M105").Value = -14912.715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3102.603
$ws.Range("I31").Value = 2042.2424
$ws.Range("K31").Value = 2042.2424
$ws.Range("M31").Value = -1747.2424
$ws.Range("H34").Value = 3102.603
$ws.Range("I34").Value = 2042.2424
$ws.Range("K34").Value = 2042.2424
$ws.Range("M34").Value = -1840.2424
$ws.Range("H58").Value = 1425.2766
$ws.Range("I58").Value = 1147.7241
$ws.Range("K58").Value = 1147.7241
$ws.Range("M58").Value = -944.7240999999999
$ws.Range("H92").Value = 29874.25
$ws.Range("J92").Value = 29874.25
$ws.Range("L92").Value = 29874.25
$ws.Range("N92").Value = -34866.25
$ws.Range("H94").Value = 4182.2812
$ws.Range("I94").Value = 4076.3333
$ws.Range("J94").Value = 4245.85
$ws.Range("K94").Value = 4076.3333
$ws.Range("L94").Value = 4245.85
$ws.Range("M94").Value = -3625.3333
$ws.Range("N94").Value = -5147.85
$ws.Range("H136").Value = 1425.2766
$ws.Range("I136").Value = 1147.7241
$ws.Range("K136").Value = 3443.1723
$ws.Range("M136").Value = -893.1722999999997
$ws.Range("H141").Value = 250071.08
$ws.Range("J141").Value = 250071.08
$ws.Range("L141").Value = 250071.08
$ws.Range("N141").Value = -260431.08

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 42111.11
$ws.Range("J37").Value = 42111.11
$ws.Range("L37").Value = 126333.33
$ws.Range("N37").Value = -126557.33
$ws.Range("H69").Value = 1774.2222
$ws.Range("I69").Value = 637.3333
$ws.Range("K69").Value = 1911.9999
$ws.Range("M69").Value = -1100.9999
$ws.Range("H72").Value = 1774.2222
$ws.Range("I72").Value = 637.3333
$ws.Range("K72").Value = 5735.9997
$ws.Range("M72").Value = -1679.9997
$ws.Range("H80").Value = 2553.7
$ws.Range("I80").Value = 2350
$ws.Range("J80").Value = 2576.3333
$ws.Range("K80").Value = 7050
$ws.Range("L80").Value = 7728.999899999999
$ws.Range("M80").Value = -6114
$ws.Range("N80").Value = -9600.999899999999
$ws.Range("H83").Value = 2553.7
$ws.Range("I83").Value = 2350
$ws.Range("J83").Value = 2576.3333
$ws.Range("K83").Value = 21150
$ws.Range("L83").Value = 23186.9997
$ws.Range("M83").Value = -16470
$ws.Range("N83").Value = -32546.9997
$ws.Range("H122").Value = 596.3077
$ws.Range("J122").Value = 610
$ws.Range("L122").Value = 5490
$ws.Range("N122").Value = -10390
$ws.Range("H131").Value = 1786686.8
$ws.Range("I131").Value = 7143496.5
$ws.Range("J131").Value = 1083.5952
$ws.Range("K131").Value = 21430489.5
$ws.Range("L131").Value = 3250.7856
$ws.Range("M131").Value = -21425449.5
$ws.Range("N131").Value = -13330.7856

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 378115.72
$ws.Range("I102").Value = 771373.9
$ws.Range("K102").Value = 771373.9
$ws.Range("M102").Value = -769751.9
$ws.Range("H126").Value = 4799.452
$ws.Range("I126").Value = 7197.8335
$ws.Range("J126").Value = 3000.6667
$ws.Range("K126").Value = 21593.5005
$ws.Range("L126").Value = 9002.000100000001
$ws.Range("M126").Value = -19123.5005
$ws.Range("N126").Value = -13942.0001
$ws.Range("H132").Value = 3211
$ws.Range("I132").Value = 4136.357
$ws.Range("J132").Value = 2529.158
$ws.Range("K132").Value = 12409.071
$ws.Range("L132").Value = 7587.474
$ws.Range("M132").Value = -9879.071
$ws.Range("N132").Value = -12647.474

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 85856.164
$ws.Range("I7").Value = 93343.09
$ws.Range("J7").Value = 3500
$ws.Range("K7").Value = 93343.09
$ws.Range("L7").Value = 3500
$ws.Range("M7").Value = -93231.09
$ws.Range("N7").Value = -3724
$ws.Range("H122").Value = 503869.84
$ws.Range("I122").Value = 4417.909
$ws.Range("J122").Value = 1114311.1
$ws.Range("K122").Value = 13253.727
$ws.Range("L122").Value = 3342933.3
$ws.Range("M122").Value = -10803.727
$ws.Range("N122").Value = -3347833.3
$ws.Range("H126").Value = 85856.164
$ws.Range("I126").Value = 93343.09
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 280029.27
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -277559.27
$ws.Range("N126").Value = -15440
$ws.Range("H136").Value = 4504.04
$ws.Range("I136").Value = 2393.6
$ws.Range("J136").Value = 7669.7
$ws.Range("K136").Value = 7180.799999999999
$ws.Range("L136").Value = 23009.1
$ws.Range("M136").Value = -4630.799999999999
$ws.Range("N136").Value = -28109.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4385.7144
$ws.Range("I62").Value = 4385.7144
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 4385.7144
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -3761.7144
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 4385.7144
$ws.Range("I65").Value = 4385.7144
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 21928.572
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -18808.572
$ws.Range("N65").ClearContents()
$ws.Range("H81").Value = 1093.8572
$ws.Range("I81").Value = 1093.8572
$ws.Range("K81").Value = 2187.7144
$ws.Range("M81").Value = -1126.7144
$ws.Range("H84").Value = 1093.8572
$ws.Range("I84").Value = 1093.8572
$ws.Range("K84").Value = 10938.572
$ws.Range("M84").Value = -5634.572
$ws.Range("H113").Value = 892.5
$ws.Range("I113").Value = 651
$ws.Range("J113").Value = 1375.5
$ws.Range("K113").Value = 1953
$ws.Range("L113").Value = 4126.5
$ws.Range("M113").Value = 217
$ws.Range("N113").Value = -8466.5
$ws.Range("H122").Value = 907.85
$ws.Range("I122").Value = 791.58826
$ws.Range("J122").Value = 1566.6666
$ws.Range("K122").Value = 2374.76478
$ws.Range("L122").Value = 4699.9998
$ws.Range("M122").Value = 75.23522000000003
$ws.Range("N122").Value = -9599.9998

